$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename a handful of activity rows: replace underscores joining codes with '+' ---
$ws.Range("A2").Value = "T+K+RTI_1"
$ws.Range("A3").Value = "T+K+RTI_2"
$ws.Range("A9").Value = "B+E_1"
$ws.Range("A10").Value = "B+E_2"
$ws.Range("A12").Value = "A+I_1"
$ws.Range("A15").Value = "s+z_1"
$ws.Range("A16").Value = "s+z_2"

# --- Row 32 (MENDEZ): MaxGroupSize -1 ---
$ws.Range("B32").Value = -1

# --- Split the old single ASSESSMENT row (old row 33) into ASSESSMENT_1 / ASSESSMENT_2 ---
# Old row33: ASSESSMENT, Duration=0, MaxGroupSize=60, Grade=1, 3 day/time blocks
# New row33: ASSESSMENT_1, Duration=-1, MaxGroupSize=30, same schedule blocks
$ws.Range("A33").Value = "ASSESSMENT_1"
$ws.Range("B33").Value = -1
$ws.Range("C33").Value = 30

# New row34: ASSESSMENT_2, Duration=-1, MaxGroupSize=30, same schedule blocks as ASSESSMENT_1
$ws.Range("A34").Value = "ASSESSMENT_2"
$ws.Range("B34").Value = -1
$ws.Range("C34").Value = 30
$ws.Range("D34").Value = 1
$ws.Range("E34").Value = 1
$ws.Range("F34").Value = $ws.Range("F33").Value2
$ws.Range("G34").Value = $ws.Range("G33").Value2
$ws.Range("H34").Value = $ws.Range("H33").Value2
$ws.Range("I34").Value = $ws.Range("I33").Value2
$ws.Range("J34").Value = $ws.Range("J33").Value2
$ws.Range("K34").Value = $ws.Range("K33").Value2
$ws.Range("L34").Value = $ws.Range("L33").Value2
$ws.Range("M34").Value = $ws.Range("M33").Value2

# match formatting of the day/start/end columns used in row 33 (center aligned;
# start/end columns use the custom time number format)
$ws.Range("H34").HorizontalAlignment = -4108
$ws.Range("K34").HorizontalAlignment = -4108
$timeFmt = "[$-F400]h:mm:ss\ AM/PM"
"I34", "J34", "L34", "M34" | ForEach-Object {
    $ws.Range($_).NumberFormat = $timeFmt
    $ws.Range($_).HorizontalAlignment = -4108
}

# --- Former LUNCH_1 / LUNCH_2 / LUNCH_3 rows (old rows 34-36) shift down one row ---
# and become new rows 35-37, each with Duration -1, MaxGroupSize shrunk from 30 to 25,
# and a ten minute earlier start / five minute later end.
$ws.Range("A35").Value = "LUNCH_1"
$ws.Range("B35").Value = -1
$ws.Range("C35").Value = 25
$ws.Range("D35").Value = 1
$ws.Range("E35").Value = 1
$ws.Range("F35").Value = 0.46527777777777773
$ws.Range("G35").Value = 0.52430555555555558

$ws.Range("A36").Value = "LUNCH_2"
$ws.Range("B36").Value = -1
$ws.Range("C36").Value = 25
$ws.Range("D36").Value = 1
$ws.Range("E36").Value = 2
$ws.Range("F36").Value = 0.46527777777777773
$ws.Range("G36").Value = 0.52430555555555558

$ws.Range("A37").Value = "LUNCH_3"
$ws.Range("B37").Value = -1
$ws.Range("C37").Value = 25
$ws.Range("D37").Value = 1
$ws.Range("E37").Value = 5
$ws.Range("F37").Value = 0.46527777777777773
$ws.Range("G37").Value = 0.50694444444444442
$ws.Range("E37").HorizontalAlignment = -4108
"F37", "G37" | ForEach-Object {
    $ws.Range($_).NumberFormat = $timeFmt
    $ws.Range($_).HorizontalAlignment = -4108
}

# --- Sheet view adjustments ---
$ws.Application.ActiveWindow.ScrollRow = 15
$ws.Application.ActiveWindow.Zoom = 150
$ws.Range("B25").Select()
